# Apply marksheet corrections: update correct/total mark counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row (row 11) - number of right answers column (B)
$ws.Range("B11").Value = 5

# "Total" row (row 12) - total marks scored (B) and the Corr/Total summary text (E)
$ws.Range("B12").Value = 90
$ws.Range("E12").Value = "90/140"
